# KHL stats refresh: re-run of scraper advancing "as_of" from 2025-10-31 to
# 2025-11-01 — appends the day's four new matches to Matches_SOG and rolls
# the date-stamped aggregate sheets (Shots_HA / Shots_Summary / Meta_ext)
# forward to match.

$wb = $excel.ActiveWorkbook

$matches   = $wb.Worksheets.Item("Matches_SOG")
$shotsHA   = $wb.Worksheets.Item("Shots_HA")
$shotsSum  = $wb.Worksheets.Item("Shots_Summary")
$metaExt   = $wb.Worksheets.Item("Meta_ext")

# ---------------------------------------------------------------------
# 1) Matches_SOG: append the 2025-11-01 matches as rows 418-421.
#    uid/date_utc/home/away/source are text columns in this sheet (the
#    uid values look numeric but are stored as text), so each gets a
#    quote-prefixed Formula assignment to force text instead of Excel's
#    auto-number coercion; sog_home/sog_away are genuine numbers.
# ---------------------------------------------------------------------
function Set-TextCell($ws, $row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Formula = "'" + $text
    # Quote-prefixing (forcing text for a numeric-looking / date-looking
    # string) otherwise stamps the cell with a "quote prefix" style; the
    # source sheets keep these columns on the default/no style, so reset it.
    $cell.Style = "Normal"
}

$newMatches = @(
    @{ Row = 418; Uid = "897718"; Date = "2025-11-01T16:30:00"; Home = "Авангард";   Away = "Нефтехимик"; Sog = 21; SogA = 31 },
    @{ Row = 419; Uid = "897716"; Date = "2025-11-01T19:00:00"; Home = "Ак Барс";    Away = "Лада";       Sog = 41; SogA = 23 },
    @{ Row = 420; Uid = "897717"; Date = "2025-11-01T19:00:00"; Home = "Северсталь"; Away = "Драконы";    Sog = 28; SogA = 23 },
    @{ Row = 421; Uid = "897715"; Date = "2025-11-01T19:30:00"; Home = "Динамо М";   Away = "ЦСКА";       Sog = 24; SogA = 19 }
)

foreach ($m in $newMatches) {
    Set-TextCell $matches $m.Row 1 $m.Uid
    Set-TextCell $matches $m.Row 2 $m.Date
    Set-TextCell $matches $m.Row 3 $m.Home
    Set-TextCell $matches $m.Row 4 $m.Away
    $matches.Cells.Item($m.Row, 5).Value = $m.Sog
    $matches.Cells.Item($m.Row, 6).Value = $m.SogA
    Set-TextCell $matches $m.Row 7 "khl_text"
}

# ---------------------------------------------------------------------
# 2) Shots_HA: as_of_utc rolls forward on every row; several teams also
#    picked up a home or away game today so their GP/OG totals move.
# ---------------------------------------------------------------------
$shotsHAUpdates = @{
    2  = @{ E = 21; G = 687;  H = 590; I = 32.7; J = 28.1 }
    3  = @{}
    4  = @{}
    5  = @{ E = 23; G = 779;  H = 587; I = 33.9; J = 25.5 }
    6  = @{}
    7  = @{}
    8  = @{ E = 15; G = 482;  H = 415; I = 32.1; J = 27.7 }
    9  = @{}
    10 = @{ F = 18; K = 503; L = 654; M = 27.9; N = 36.3 }
    11 = @{ F = 17; K = 457; L = 623; M = 26.9; N = 36.6 }
    12 = @{}
    13 = @{}
    14 = @{ F = 17; K = 465; L = 643; M = 27.4; N = 37.8 }
    15 = @{}
    16 = @{}
    17 = @{ E = 16; G = 469;  H = 363; I = 29.3 }
    18 = @{}
    19 = @{}
    20 = @{}
    21 = @{}
    22 = @{}
    23 = @{ F = 20; K = 497; L = 574; M = 24.9; N = 28.7 }
}

foreach ($row in $shotsHAUpdates.Keys) {
    Set-TextCell $shotsHA $row 4 "2025-11-01T19:30:00Z"
    foreach ($entry in $shotsHAUpdates[$row].GetEnumerator()) {
        $col = switch ($entry.Key) {
            "E" { 5 }; "F" { 6 }; "G" { 7 }; "H" { 8 }; "I" { 9 }
            "J" { 10 }; "K" { 11 }; "L" { 12 }; "M" { 13 }; "N" { 14 }
        }
        $shotsHA.Cells.Item($row, $col).Value = $entry.Value
    }
}

# ---------------------------------------------------------------------
# 3) Shots_Summary: same as_of_utc roll, same teams' SOG/SOGA totals move.
# ---------------------------------------------------------------------
$shotsSumUpdates = @{
    2  = @{ E = 37; F = 1259; G = 1062; H = 34;   I = 28.7 }
    3  = @{}
    4  = @{}
    5  = @{ E = 41; F = 1383; G = 1121; H = 33.7; I = 27.3 }
    6  = @{}
    7  = @{}
    8  = @{ E = 35; F = 1042; G = 1065; H = 29.8; I = 30.4 }
    9  = @{}
    10 = @{ E = 37; F = 1035; G = 1316; H = 28;   I = 35.6 }
    11 = @{ E = 39; F = 1055; G = 1399; H = 27.1; I = 35.9 }
    12 = @{}
    13 = @{}
    14 = @{ E = 40; F = 1189; G = 1413;           I = 35.3 }
    15 = @{}
    16 = @{}
    17 = @{ E = 38; F = 1191; G = 943;  H = 31.3; I = 24.8 }
    18 = @{}
    19 = @{}
    20 = @{}
    21 = @{}
    22 = @{}
    23 = @{ E = 37; F = 880;  G = 1073; H = 23.8; I = 29 }
}

foreach ($row in $shotsSumUpdates.Keys) {
    Set-TextCell $shotsSum $row 4 "2025-11-01T19:30:00Z"
    foreach ($entry in $shotsSumUpdates[$row].GetEnumerator()) {
        $col = switch ($entry.Key) {
            "E" { 5 }; "F" { 6 }; "G" { 7 }; "H" { 8 }; "I" { 9 }
        }
        $shotsSum.Cells.Item($row, $col).Value = $entry.Value
    }
}

# ---------------------------------------------------------------------
# 4) Meta_ext: bump as_of_utc and build_version.
# ---------------------------------------------------------------------
Set-TextCell $metaExt 2 2 "2025-11-01T19:30:00Z"
$metaExt.Cells.Item(2, 4).Value = 32
